# Generate Report for Handoff
#
# - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
# - "Latest HO Xliff Generate Date" / "Latest Handback DateTime" timestamp
#   2016-09-02 01:09:49 -> 2016-09-02 01:10:49 (Overview!G2, de-de!H2)
# - zh-cn "Latest Handoff Datetime" timestamp
#   2016-09-02 01:09:45 -> 2016-09-02 01:10:45 (zh-cn!H2)
# - Narrow the (now shorter) status columns: Overview E:F and the "Status"
#   column (C) on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$oldHandback = "2016-09-02 01:09:49"
$newHandback = "2016-09-02 01:10:49"

$oldHandoff = "2016-09-02 01:09:45"
$newHandoff = "2016-09-02 01:10:45"

# Target stored OOXML column width is 17.2159881591797 characters. Excel's
# ColumnWidth setter always quantizes to whole screen pixels (width_px =
# round(ColumnWidth * MaxDigitWidth), then stored_width = (width_px + 5) /
# MaxDigitWidth for the workbook's default Calibri 11 font, MaxDigitWidth =
# 6px) before it is written back to the sheet, so the exact fractional
# value can't be round-tripped through the object model. 16.333333333333332
# is the ColumnWidth input that quantizes to the closest reachable stored
# width (17.166666666666668, i.e. 98px) to the target.
$newColWidth = 16.333333333333332

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHandback

$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newHandoff

$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newHandback

$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
